$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.55
$ws.Range("H2").Value = 4
$ws.Range("I2").Value = 6
$ws.Range("J2").Value = 2.2
$ws.Range("L2").Value = 6.5
$ws.Range("O2").Value = 1.36
$ws.Range("P2").Value = 3
$ws.Range("S2").Value = 1.5
$ws.Range("T2").Value = 2.5
$ws.Range("U2").Value = 2.25
$ws.Range("V2").Value = 1.57
$ws.Range("X2").Value = 6
$ws.Range("Z2").Value = 11
$ws.Range("AD2").Value = 8
$ws.Range("AE2").Value = 23
$ws.Range("AH2").Value = 12
$ws.Range("AI2").Value = 29
$ws.Range("AJ2").Value = 19
$ws.Range("AK2").Value = 67
$ws.Range("AL2").Value = 51
$ws.Range("AN2").Value = 3.25
$ws.Range("AO2").Value = 8
$ws.Range("AQ2").Value = 26
$ws.Range("AT2").Value = 2.5
$ws.Range("AU2").Value = 10
$ws.Range("AV2").Value = 81
$ws.Range("AW2").Value = 7.5
$ws.Range("AZ2").Value = 151
$ws.Range("BA2").Value = 201
$ws.Range("G3").Value = 1.48
$ws.Range("H3").Value = 4
$ws.Range("I3").Value = 7
$ws.Range("J3").Value = 2.1
$ws.Range("L3").Value = 7.5
$ws.Range("O3").Value = 1.36
$ws.Range("P3").Value = 3
$ws.Range("S3").Value = 1.44
$ws.Range("T3").Value = 2.63
$ws.Range("U3").Value = 2.38
$ws.Range("V3").Value = 1.53
$ws.Range("W3").Value = 5.5
$ws.Range("X3").Value = 6
$ws.Range("Y3").Value = 9
$ws.Range("Z3").Value = 9.5
$ws.Range("AD3").Value = 8
$ws.Range("AE3").Value = 23
$ws.Range("AK3").Value = 81
$ws.Range("AL3").Value = 51
$ws.Range("AS3").Value = 201
$ws.Range("AT3").Value = 2.63
$ws.Range("AU3").Value = 10
$ws.Range("AW3").Value = 8
$ws.Range("AY3").Value = 41
$ws.Range("G4").Value = 2.8
$ws.Range("H4").Value = 3.1
$ws.Range("I4").Value = 2.45
$ws.Range("K4").Value = 2.1
$ws.Range("L4").Value = 3
$ws.Range("N4").Value = 8.1
$ws.Range("O4").Value = 1.29
$ws.Range("P4").Value = 2.95
$ws.Range("Q4").Value = 1.93
$ws.Range("R4").Value = 1.78
$ws.Range("U4").Value = 1.65
$ws.Range("V4").Value = 1.98
$ws.Range("W4").Value = 9.25
$ws.Range("AB4").Value = 29
$ws.Range("AC4").Value = 9.25
$ws.Range("AE4").Value = 12.5
$ws.Range("AF4").Value = 55
$ws.Range("AG4").Value = 400
$ws.Range("AH4").Value = 8
$ws.Range("AI4").Value = 12.5
$ws.Range("AJ4").Value = 9.25
$ws.Range("AM4").Value = 29
$ws.Range("AN4").Value = 4.8
$ws.Range("AP4").Value = 19
$ws.Range("AT4").Value = 2.65
$ws.Range("AU4").Value = 6.4
$ws.Range("AX4").Value = 12.5
$ws.Range("AY4").Value = 18.5
$ws.Range("AZ4").Value = 50
$ws.Range("BA4").Value = 75
$ws.Range("BB4").Value = 200
